$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily data rows for LamaMocogno covering through 2021-12-08 (aggiornamento fino a 8/12)
# Columns: row, date-serial (A), nuovi pos. (B), somma mobile 7gg. (C), somma mobile 7gg. per 100mila abitanti (D)
$data = @(
  @(386, 44460, 0, 2, 74.93443237167479),
  @(387, 44461, 0, 2, 74.93443237167479),
  @(388, 44462, 0, 2, 74.93443237167479),
  @(389, 44463, 1, 1, 37.46721618583739),
  @(390, 44464, 0, 1, 37.46721618583739),
  @(391, 44465, 0, 1, 37.46721618583739),
  @(392, 44466, 0, 1, 37.46721618583739),
  @(393, 44467, 1, 2, 74.93443237167479),
  @(394, 44468, 0, 2, 74.93443237167479),
  @(395, 44469, 0, 2, 74.93443237167479),
  @(396, 44470, 1, 2, 74.93443237167479),
  @(397, 44471, 0, 2, 74.93443237167479),
  @(398, 44472, 0, 2, 74.93443237167479),
  @(399, 44473, 0, 2, 74.93443237167479),
  @(400, 44474, 0, 1, 37.46721618583739),
  @(401, 44475, 0, 1, 37.46721618583739),
  @(402, 44476, 1, 2, 74.93443237167479),
  @(403, 44477, 0, 1, 37.46721618583739),
  @(404, 44478, 0, 1, 37.46721618583739),
  @(405, 44479, 0, 1, 37.46721618583739),
  @(406, 44480, 0, 1, 37.46721618583739),
  @(407, 44481, 0, 1, 37.46721618583739),
  @(408, 44482, 0, 1, 37.46721618583739),
  @(409, 44483, 0, 0, 0),
  @(410, 44484, 0, 0, 0),
  @(411, 44485, 0, 0, 0),
  @(412, 44486, 0, 0, 0),
  @(413, 44487, 0, 0, 0),
  @(414, 44488, 0, 0, 0),
  @(415, 44489, 0, 0, 0),
  @(416, 44490, 0, 0, 0),
  @(417, 44491, 0, 0, 0),
  @(418, 44492, 0, 0, 0),
  @(419, 44493, 0, 0, 0),
  @(420, 44494, 0, 0, 0),
  @(421, 44495, 0, 0, 0),
  @(422, 44496, 0, 0, 0),
  @(423, 44497, 0, 0, 0),
  @(424, 44498, 0, 0, 0),
  @(425, 44499, 0, 0, 0),
  @(426, 44500, 0, 0, 0),
  @(427, 44501, 0, 0, 0),
  @(428, 44502, 0, 0, 0),
  @(429, 44503, 0, 0, 0),
  @(430, 44504, 0, 0, 0),
  @(431, 44505, 0, 0, 0),
  @(432, 44506, 1, 1, 37.46721618583739),
  @(433, 44507, 0, 1, 37.46721618583739),
  @(434, 44508, 0, 1, 37.46721618583739),
  @(435, 44509, 0, 1, 37.46721618583739),
  @(436, 44510, 0, 1, 37.46721618583739),
  @(437, 44511, 0, 1, 37.46721618583739),
  @(438, 44512, 0, 1, 37.46721618583739),
  @(439, 44513, 0, 0, 0),
  @(440, 44514, 0, 0, 0),
  @(441, 44515, 0, 0, 0),
  @(442, 44516, 2, 2, 74.93443237167479),
  @(443, 44517, 0, 2, 74.93443237167479),
  @(444, 44518, 0, 2, 74.93443237167479),
  @(445, 44519, 0, 2, 74.93443237167479),
  @(446, 44520, 0, 2, 74.93443237167479),
  @(447, 44521, 0, 2, 74.93443237167479),
  @(448, 44522, 0, 2, 74.93443237167479),
  @(449, 44523, 0, 0, 0),
  @(450, 44524, 2, 2, 74.93443237167479),
  @(451, 44525, 0, 2, 74.93443237167479),
  @(452, 44526, 0, 2, 74.93443237167479),
  @(453, 44527, 0, 2, 74.93443237167479),
  @(454, 44528, 0, 2, 74.93443237167479),
  @(455, 44529, 0, 2, 74.93443237167479),
  @(456, 44530, 0, 2, 74.93443237167479),
  @(457, 44531, 0, 0, 0),
  @(458, 44532, 0, 0, 0),
  @(459, 44533, 1, 1, 37.46721618583739),
  @(460, 44534, 0, 1, 37.46721618583739),
  @(461, 44535, 0, 1, 37.46721618583739),
  @(462, 44536, 0, 1, 37.46721618583739),
  @(463, 44537, 0, 1, 37.46721618583739),
  @(464, 44538, 1, 2, 74.93443237167479)
)

$lastRow = 385
foreach ($item in $data) {
    $r = $item[0]
    $dateSerial = $item[1]
    $newPos = $item[2]
    $rollingSum = $item[3]
    $rollingPer100k = $item[4]

    # Carry the existing date-column formatting (border + center/top alignment + datetime numberformat)
    # down onto the newly appended row, matching how the prior rows were styled.
    $ws.Range("A$lastRow").Copy($ws.Range("A$r"))

    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $newPos
    $ws.Cells.Item($r, 3).Value = $rollingSum
    $ws.Cells.Item($r, 4).Value = $rollingPer100k

    $lastRow = $r
}
